$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Exercitiile de la curs mai explicate" - a batch of students picked up an
# extra attendance point for "saptamana 3" (column E), bumping it from 1 to 2.
# The "Prezente" (Q) column is a SUM formula, so it recalculates automatically.
$rows = @(3, 6, 8, 12, 13, 14, 17, 18, 19, 21)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 5).Value = 2
}

# Update the cell that is selected/active when the workbook is reopened.
[void]$ws.Range("H20").Select()
